# 2025 author list / membership refresh for the Strathclyde collaboration list.
# - Add the JISCMail alias to Paul McKenna's email address
# - Add two new members (Ewan Dolier, Matthew Alderton) to the table, each with
#   a mailto: hyperlink on their email address
# - Keep the table / autofilter / dimension in sync with the new rows
# - Re-select the cell the author ended up on (K5)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Paul McKenna now also has his JISCMail alias on file -----------------
$ws.Range("E2").Value = $ws.Range("E2").Value2 + "; p.mckenna@phys.strath.ac.uk"

# --- 2. Grow the table by two rows so formatting / autofilter follow along ---
$tbl = $ws.ListObjects.Item(1)
$tbl.ListRows.Add() | Out-Null
$tbl.ListRows.Add() | Out-Null

# --- 3. New member: Ewan Dolier (row 8) --------------------------------------
$ws.Range("A8").Value = "Dr."
$ws.Range("B8").Value = "Ewan"
$ws.Range("C8").Value = "Dolier"
$ws.Range("D8").Value = "E."
$ws.Range("E8").Value = "ewan.dolier.2015@uni.strath.ac.uk"
$ws.Range("F8").Value = "E.Dolier"
$ws.Range("G8").Value = "Strathclyde"
$ws.Range("H8").Value = "Department of Physics, SUPA, University of Strathclyde, Glasgow G4 0NG, UK"
$ws.Range("I8").Value = 0

# --- 4. New member: Matthew Alderton (row 9) ---------------------------------
$ws.Range("A9").Value = "Dr."
$ws.Range("B9").Value = "Matthew"
$ws.Range("C9").Value = "Alderton"
$ws.Range("D9").Value = "M."
$ws.Range("E9").Value = "matthew.alderton@strath.ac.uk "
$ws.Range("F9").Value = "M.Alderton"
$ws.Range("G9").Value = "Strathclyde"
$ws.Range("H9").Value = "Department of Physics, SUPA, University of Strathclyde, Glasgow G4 0NG, UK"
$ws.Range("I9").Value = 0

# --- 5. Hyperlink the two new email addresses (mailto:) ----------------------
$ws.Hyperlinks.Add($ws.Range("E8"), "mailto:ewan.dolier.2015@uni.strath.ac.uk")
$ws.Hyperlinks.Add($ws.Range("E9"), "mailto:matthew.alderton@strath.ac.uk")

# --- 6. Widen the email-address column to fit the longer entries -------------
$ws.Columns.Item(5).ColumnWidth = 49.25

# --- 7. Leave the selection where the author left it --------------------------
$ws.Range("K5").Select() | Out-Null
